{"js": "// Shorten the opening \"meta description\" paragraph so it reads under 160\n// characters, per the commit message. The only visible content change in\n// the target revision is the wording of the very first paragraph:\n//\n//   \"A collection of projects I created during my sophomore year at RISD\n//    as an industrial design major. These include the ...\"\n// becomes\n//   \"Projects I created during my RISD sophomore year as an industrial\n//    design major. These include the ...\"\n//\n// (Everything else in the source diff is Word re-splitting runs to wrap\n// w:proofErr spell/grammar markers around unchanged words such as\n// \"french\"/\"Klann\"/\"Balani\" \u2014 the visible text is identical, so there is\n// nothing further to edit.)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\nconst OLD_PREFIX =\n  \"A collection of projects I created during my sophomore year at RISD as an industrial design major. These include the \";\nconst NEW_PREFIX =\n  \"Projects I created during my RISD sophomore year as an industrial design major. These include the \";\n\nconst searchResults = firstParagraph.search(OLD_PREFIX, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(NEW_PREFIX, Word.InsertLocation.replace);\n} else {\n  // Fallback: in case the exact prefix isn't found as a single contiguous\n  // match (e.g. already edited), rebuild the paragraph text explicitly.\n  firstParagraph.load(\"text\");\n  await context.sync();\n  if (firstParagraph.text.indexOf(OLD_PREFIX) === 0) {\n    const rest = firstParagraph.text.substring(OLD_PREFIX.length);\n    firstParagraph.insertText(NEW_PREFIX + rest, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Shorten the opening \"meta description\" paragraph so it reads under 160\n# characters, per the commit message. The only visible content change in\n# the target revision is the wording of the very first paragraph:\n#\n#   \"A collection of projects I created during my sophomore year at RISD\n#    as an industrial design major. These include the ...\"\n# becomes\n#   \"Projects I created during my RISD sophomore year as an industrial\n#    design major. These include the ...\"\n#\n# (Everything else in the source diff is Word re-splitting runs to wrap\n# w:proofErr spell/grammar markers around unchanged words such as\n# \"french\"/\"Klann\"/\"Balani\" -- the visible text is identical, so there is\n# nothing further to edit.)\n\n$d = $word.ActiveDocument\n\n$oldPrefix = \"A collection of projects I created during my sophomore year at RISD as an industrial design major. These include the \"\n$newPrefix = \"Projects I created during my RISD sophomore year as an industrial design major. These include the \"\n\n# Scope the search to the first paragraph only, so the similarly-worded\n# paragraph further down (\"...Most of these remain mere explorations...\")\n# is left untouched.\n$firstParagraph = $d.Paragraphs.Item(1).Range\n$searchRange = $d.Range($firstParagraph.Start, $firstParagraph.End)\n\n$searchRange.Find.ClearFormatting()\n$searchRange.Find.Replacement.ClearFormatting()\n$found = $searchRange.Find.Execute(\n    $oldPrefix,\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    $newPrefix, 2\n)\n\nif (-not $found) {\n    # Fallback in case the exact prefix can't be located as one contiguous\n    # match (e.g. the document was already edited) -- rebuild the text\n    # directly from the paragraph's current content.\n    $text = $firstParagraph.Text\n    if ($text.StartsWith($oldPrefix)) {\n        $rest = $text.Substring($oldPrefix.Length)\n        $firstParagraph.Text = $newPrefix + $rest\n    }\n}\n"}
